$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '27.052.78'
$ws.Range("E2").NumberFormat = '@'
$ws.Range("E2").Value = '  -0.48%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '1.829.10'
$ws.Range("E3").NumberFormat = '@'
$ws.Range("E3").Value = '  +0.03%  '

$ws.Range("D4").NumberFormat = '@'
$ws.Range("D4").Value = '1.008'
$ws.Range("E4").NumberFormat = '@'
$ws.Range("E4").Value = '  -0.31%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '312.30'
$ws.Range("E5").NumberFormat = '@'
$ws.Range("E5").Value = '  -0.44%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '1.007'
$ws.Range("E6").NumberFormat = '@'
$ws.Range("E6").Value = '  -0.47%  '

$ws.Range("E7").NumberFormat = '@'
$ws.Range("E7").Value = '  -0.88%  '

$ws.Range("D8").NumberFormat = '@'
$ws.Range("D8").Value = '0.3700'
$ws.Range("E8").NumberFormat = '@'
$ws.Range("E8").Value = '  +1.59%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.07341'
$ws.Range("E9").NumberFormat = '@'
$ws.Range("E9").Value = '  +0.12%  '

$ws.Range("D10").NumberFormat = '@'
$ws.Range("D10").Value = '0.8707'
$ws.Range("E10").NumberFormat = '@'
$ws.Range("E10").Value = '  -0.67%  '

$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '0.07920'
$ws.Range("E11").NumberFormat = '@'
$ws.Range("E11").Value = '  +3.49%  '

$ws.Range("D12").NumberFormat = '@'
$ws.Range("D12").Value = '19.79'
$ws.Range("E12").NumberFormat = '@'
$ws.Range("E12").Value = '  -2.16%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '1.849.47'
$ws.Range("E13").NumberFormat = '@'
$ws.Range("E13").Value = '  -0.16%  '

$ws.Range("D14").NumberFormat = '@'
$ws.Range("D14").Value = '5.342'
$ws.Range("E14").NumberFormat = '@'
$ws.Range("E14").Value = '  -0.40%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '6.551'
$ws.Range("E15").NumberFormat = '@'
$ws.Range("E15").Value = '  +0.96%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '91.72'
$ws.Range("E16").NumberFormat = '@'
$ws.Range("E16").Value = '  -1.20%  '

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '1.009'
$ws.Range("E17").NumberFormat = '@'
$ws.Range("E17").Value = '  +0.00%  '

$ws.Range("D18").NumberFormat = '@'
$ws.Range("D18").Value = '0.000008867'
$ws.Range("E18").NumberFormat = '@'
$ws.Range("E18").Value = '  +2.10%  '

$ws.Range("E19").NumberFormat = '@'
$ws.Range("E19").Value = '  -0.50%  '

$ws.Range("E20").NumberFormat = '@'
$ws.Range("E20").Value = '  +1.06%  '

$ws.Range("D21").NumberFormat = '@'
$ws.Range("D21").Value = '26.932.62'
$ws.Range("E21").NumberFormat = '@'
$ws.Range("E21").Value = '  -2.32%  '

$ws.Range("D22").NumberFormat = '@'
$ws.Range("D22").Value = '5.120'
$ws.Range("E22").NumberFormat = '@'
$ws.Range("E22").Value = '  -2.16%  '

$ws.Range("E23").NumberFormat = '@'
$ws.Range("E23").Value = '  -0.21%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '2.002.89'
$ws.Range("E24").NumberFormat = '@'
$ws.Range("E24").Value = '  -4.53%  '

$ws.Range("D25").NumberFormat = '@'
$ws.Range("D25").Value = '152.48'
$ws.Range("E25").NumberFormat = '@'
$ws.Range("E25").Value = '  +0.80%  '

$ws.Range("D26").NumberFormat = '@'
$ws.Range("D26").Value = '1.841'
$ws.Range("E26").NumberFormat = '@'
$ws.Range("E26").Value = '  -2.28%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '18.50'
$ws.Range("E27").NumberFormat = '@'
$ws.Range("E27").Value = '  +0.25%  '

$ws.Range("E28").NumberFormat = '@'
$ws.Range("E28").Value = '  -1.30%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '5.104'
$ws.Range("E29").NumberFormat = '@'
$ws.Range("E29").Value = '  -0.59%  '

$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '115.08'
$ws.Range("E30").NumberFormat = '@'
$ws.Range("E30").Value = '  -1.10%  '

$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '0.08862'
$ws.Range("E31").NumberFormat = '@'
$ws.Range("E31").Value = '  -0.66%  '

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '2.974'
$ws.Range("E32").NumberFormat = '@'
$ws.Range("E32").Value = '  +0.40%  '

$ws.Range("D33").NumberFormat = '@'
$ws.Range("D33").Value = '0.7312'
$ws.Range("E33").NumberFormat = '@'
$ws.Range("E33").Value = '  -0.93%  '

$ws.Range("D34").NumberFormat = '@'
$ws.Range("D34").Value = '4.434'
$ws.Range("E34").NumberFormat = '@'
$ws.Range("E34").Value = '  -1.10%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '1.134'
$ws.Range("E35").NumberFormat = '@'
$ws.Range("E35").Value = '  -1.88%  '

$ws.Range("D36").NumberFormat = '@'
$ws.Range("D36").Value = '2.448'
$ws.Range("E36").NumberFormat = '@'
$ws.Range("E36").Value = '  -3.93%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '1.074'
$ws.Range("E37").NumberFormat = '@'
$ws.Range("E37").Value = '  -1.36%  '

$ws.Range("D38").NumberFormat = '@'
$ws.Range("D38").Value = '0.05242'
$ws.Range("E38").NumberFormat = '@'
$ws.Range("E38").Value = '  -0.38%  '

$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '0.01935'
$ws.Range("E39").NumberFormat = '@'
$ws.Range("E39").Value = '  +0.58%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '2.947'
$ws.Range("E40").NumberFormat = '@'
$ws.Range("E40").Value = '  +0.40%  '

$ws.Range("D41").NumberFormat = '@'
$ws.Range("D41").Value = '7.158'
$ws.Range("E41").NumberFormat = '@'
$ws.Range("E41").Value = '  -1.66%  '

$ws.Range("D42").NumberFormat = '@'
$ws.Range("D42").Value = '0.5143'
$ws.Range("E42").NumberFormat = '@'
$ws.Range("E42").Value = '  -1.95%  '

$ws.Range("E43").NumberFormat = '@'
$ws.Range("E43").Value = '  -0.16%  '

$ws.Range("D44").NumberFormat = '@'
$ws.Range("D44").Value = '0.8579'
$ws.Range("E44").NumberFormat = '@'
$ws.Range("E44").Value = '  -15.41%  '

$ws.Range("D45").NumberFormat = '@'
$ws.Range("D45").Value = '8.230'
$ws.Range("E45").NumberFormat = '@'
$ws.Range("E45").Value = '  -1.16%  '

$ws.Range("E46").NumberFormat = '@'
$ws.Range("E46").Value = '  -0.86%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = '@'
$ws.Range("D47").Value = '10.26'
$ws.Range("E47").NumberFormat = '@'
$ws.Range("E47").Value = '  -0.16%  '

$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D48").NumberFormat = '@'
$ws.Range("D48").Value = '1.007'
$ws.Range("E48").NumberFormat = '@'
$ws.Range("E48").Value = '  -0.61%  '

$ws.Range("D49").NumberFormat = '@'
$ws.Range("D49").Value = '102.27'
$ws.Range("E49").NumberFormat = '@'
$ws.Range("E49").Value = '  -1.40%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '1.625'
$ws.Range("E50").NumberFormat = '@'
$ws.Range("E50").Value = '  -1.01%  '

$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '0.06219'
$ws.Range("E51").NumberFormat = '@'
$ws.Range("E51").Value = '  -0.96%  '
